$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column A values (rows 2-11) from 1 to 2
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = 2
}

# Update the active cell selection to C12
$ws.Range("C12").Select()
